# Apply "enhanced course legend with ltpsc structure" change:
# Session duration moves from 2 hours (120 min) to 3 hours (180 min),
# which shifts every "09:00 - 11:00" slot to "09:00 - 12:00" and every
# "14:00 - 16:00" slot to "14:00 - 17:00" across the workbook, and the
# downstream summary totals that depend on duration are recomputed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Exam_Schedule
#   D: duration ("2 hours" -> "3 hours")
#   E: duration_minutes (120 -> 180)
#   K: time_slot ("09:00 - 11:00" -> "09:00 - 12:00",
#                 "14:00 - 16:00" -> "14:00 - 17:00")
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Exam_Schedule")
$lastRowSchedule = $wsSchedule.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRowSchedule; $r++) {
    $duration = $wsSchedule.Cells.Item($r, 4).Value2
    if ($duration -eq "2 hours") {
        $wsSchedule.Cells.Item($r, 4).Value = "3 hours"
    }

    $durationMinutes = $wsSchedule.Cells.Item($r, 5).Value2
    if ($durationMinutes -eq 120) {
        $wsSchedule.Cells.Item($r, 5).Value = 180
    }

    $timeSlot = $wsSchedule.Cells.Item($r, 11).Value2
    if ($timeSlot -eq "09:00 - 11:00") {
        $wsSchedule.Cells.Item($r, 11).Value = "09:00 - 12:00"
    } elseif ($timeSlot -eq "14:00 - 16:00") {
        $wsSchedule.Cells.Item($r, 11).Value = "14:00 - 17:00"
    }
}

# ---------------------------------------------------------------------
# Sheet 2: Exam_Classrooms
#   E: Time Slot ("09:00 - 11:00" -> "09:00 - 12:00",
#                 "14:00 - 16:00" -> "14:00 - 17:00")
#   I: Duration ("2 hours" -> "3 hours")
# ---------------------------------------------------------------------
$wsClassrooms = $wb.Worksheets.Item("Exam_Classrooms")
$lastRowClassrooms = $wsClassrooms.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRowClassrooms; $r++) {
    $timeSlot = $wsClassrooms.Cells.Item($r, 5).Value2
    if ($timeSlot -eq "09:00 - 11:00") {
        $wsClassrooms.Cells.Item($r, 5).Value = "09:00 - 12:00"
    } elseif ($timeSlot -eq "14:00 - 16:00") {
        $wsClassrooms.Cells.Item($r, 5).Value = "14:00 - 17:00"
    }

    $duration = $wsClassrooms.Cells.Item($r, 9).Value2
    if ($duration -eq "2 hours") {
        $wsClassrooms.Cells.Item($r, 9).Value = "3 hours"
    }
}

# ---------------------------------------------------------------------
# Sheet 3: Configuration
#   B3: Session Duration (minutes) (120 -> 180)
# ---------------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("Configuration")
$wsConfig.Range("B3").Value = 180

# ---------------------------------------------------------------------
# Sheet 5: Department_Summary
#   C: Total Duration (min) = Number of Exams * 180 (was * 120)
#   E: Total Duration (hours) = Total Duration (min) / 60
# ---------------------------------------------------------------------
$wsDept = $wb.Worksheets.Item("Department_Summary")
$lastRowDept = $wsDept.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRowDept; $r++) {
    $numExams = $wsDept.Cells.Item($r, 2).Value2
    $newTotalMinutes = $numExams * 180
    $wsDept.Cells.Item($r, 3).Value = $newTotalMinutes
    $wsDept.Cells.Item($r, 5).Value = $newTotalMinutes / 60
}
